$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-07 Friday" "2025-03-08 Saturday"

Replace-Text "427÷3=" "875÷7="
Replace-Text "812÷3=" "244÷2="
Replace-Text "733÷4=" "740÷2="
Replace-Text "359÷2=" "973÷4="
Replace-Text "540÷8=" "222÷9="
Replace-Text "820÷4=" "838÷5="
Replace-Text "541÷3=" "770÷3="
Replace-Text "846÷7=" "221÷7="
Replace-Text "742÷6=" "874÷8="
Replace-Text "835÷8=" "982÷8="
Replace-Text "416÷8=" "387÷2="
Replace-Text "416÷3=" "597÷8="
Replace-Text "634÷9=" "904÷8="
Replace-Text "971÷4=" "163÷6="
Replace-Text "418÷5=" "131÷9="
Replace-Text "938÷5=" "176÷8="
Replace-Text "897÷6=" "988÷4="
Replace-Text "313÷7=" "346÷2="
Replace-Text "517÷9=" "509÷5="
Replace-Text "418÷3=" "878÷5="
Replace-Text "587÷7=" "994÷5="
Replace-Text "514÷3=" "166÷6="
Replace-Text "842÷2=" "138÷8="
Replace-Text "186÷8=" "137÷4="
Replace-Text "606÷3=" "690÷9="
